# EXPORT aggiunta gestione mancanti e aggiunta in DDT
# Update "LANCIO" (batch id) and "PAIA DA PRODURRE" (pairs to produce),
# then recompute the TOTALE column (F) for every consumption row as
# CONS/PA (E) * PAIA DA PRODURRE (B3), rounded to 2 decimals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header values
$ws.Range("B2").Value2 = 6664
$ws.Range("B3").Value2 = 125

$paia = $ws.Range("B3").Value2

# Rows that contain a CONS/PA (E) / TOTALE (F) pair needing recalculation
$rows = @(7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38)

foreach ($r in $rows) {
    $cons = $ws.Cells.Item($r, 5).Value2
    if ($cons -ne $null) {
        $total = [Math]::Round($cons * $paia, 2)
        $ws.Cells.Item($r, 6).Value2 = $total
    }
}
